$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-05-06 Monday" "2024-05-07 Tuesday"

Replace-Text "147÷7=21, 0" "377÷7=53, 6"
Replace-Text "209÷7=29, 6" "400÷7=57, 1"
Replace-Text "188÷5=37, 3" "467÷8=58, 3"
Replace-Text "886÷3=295, 1" "109÷6=18, 1"
Replace-Text "794÷3=264, 2" "413÷6=68, 5"

Replace-Text "982÷7=140, 2" "688÷9=76, 4"
Replace-Text "709÷8=88, 5" "661÷5=132, 1"
Replace-Text "291÷4=72, 3" "826÷5=165, 1"
Replace-Text "230÷2=115, 0" "971÷2=485, 1"
Replace-Text "734÷6=122, 2" "151÷8=18, 7"

Replace-Text "871÷8=108, 7" "121÷3=40, 1"
Replace-Text "910÷4=227, 2" "862÷6=143, 4"
Replace-Text "357÷5=71, 2" "128÷5=25, 3"
Replace-Text "187÷6=31, 1" "484÷3=161, 1"
Replace-Text "232÷6=38, 4" "639÷8=79, 7"

Replace-Text "924÷9=102, 6" "899÷3=299, 2"
Replace-Text "939÷7=134, 1" "695÷7=99, 2"
Replace-Text "872÷8=109, 0" "274÷3=91, 1"
Replace-Text "588÷3=196, 0" "711÷5=142, 1"
Replace-Text "712÷6=118, 4" "186÷8=23, 2"

Replace-Text "241÷7=34, 3" "623÷5=124, 3"
Replace-Text "741÷3=247, 0" "413÷7=59, 0"
Replace-Text "403÷9=44, 7" "250÷3=83, 1"
Replace-Text "488÷8=61, 0" "843÷5=168, 3"
Replace-Text "929÷3=309, 2" "827÷3=275, 2"

$d.Save()
